$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-03 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-04 Monday", 2)

$d.Content.Find.Execute("822×9=7398", $true, $false, $false, $false, $false, $true, 1, $false, "389×7=2723", 2)
$d.Content.Find.Execute("864×4=3456", $true, $false, $false, $false, $false, $true, 1, $false, "266×8=2128", 2)
$d.Content.Find.Execute("829×9=7461", $true, $false, $false, $false, $false, $true, 1, $false, "384×9=3456", 2)
$d.Content.Find.Execute("425×7=2975", $true, $false, $false, $false, $false, $true, 1, $false, "163×8=1304", 2)
$d.Content.Find.Execute("241×2=482", $true, $false, $false, $false, $false, $true, 1, $false, "734×5=3670", 2)

$d.Content.Find.Execute("926×2=1852", $true, $false, $false, $false, $false, $true, 1, $false, "154×6=924", 2)
$d.Content.Find.Execute("242×3=726", $true, $false, $false, $false, $false, $true, 1, $false, "323×9=2907", 2)
$d.Content.Find.Execute("614×2=1228", $true, $false, $false, $false, $false, $true, 1, $false, "234×4=936", 2)
$d.Content.Find.Execute("273×2=546", $true, $false, $false, $false, $false, $true, 1, $false, "810×5=4050", 2)
$d.Content.Find.Execute("925×6=5550", $true, $false, $false, $false, $false, $true, 1, $false, "401×9=3609", 2)

$d.Content.Find.Execute("717×8=5736", $true, $false, $false, $false, $false, $true, 1, $false, "977×9=8793", 2)
$d.Content.Find.Execute("931×6=5586", $true, $false, $false, $false, $false, $true, 1, $false, "365×5=1825", 2)
$d.Content.Find.Execute("497×7=3479", $true, $false, $false, $false, $false, $true, 1, $false, "693×6=4158", 2)
$d.Content.Find.Execute("538×7=3766", $true, $false, $false, $false, $false, $true, 1, $false, "636×2=1272", 2)
$d.Content.Find.Execute("862×4=3448", $true, $false, $false, $false, $false, $true, 1, $false, "280×3=840", 2)

$d.Content.Find.Execute("716×9=6444", $true, $false, $false, $false, $false, $true, 1, $false, "458×6=2748", 2)
$d.Content.Find.Execute("986×6=5916", $true, $false, $false, $false, $false, $true, 1, $false, "133×7=931", 2)
$d.Content.Find.Execute("239×2=478", $true, $false, $false, $false, $false, $true, 1, $false, "867×5=4335", 2)
$d.Content.Find.Execute("579×2=1158", $true, $false, $false, $false, $false, $true, 1, $false, "984×9=8856", 2)
$d.Content.Find.Execute("235×9=2115", $true, $false, $false, $false, $false, $true, 1, $false, "205×8=1640", 2)

$d.Content.Find.Execute("299×3=897", $true, $false, $false, $false, $false, $true, 1, $false, "495×5=2475", 2)
$d.Content.Find.Execute("229×5=1145", $true, $false, $false, $false, $false, $true, 1, $false, "939×7=6573", 2)
$d.Content.Find.Execute("886×3=2658", $true, $false, $false, $false, $false, $true, 1, $false, "757×2=1514", 2)
$d.Content.Find.Execute("390×4=1560", $true, $false, $false, $false, $false, $true, 1, $false, "394×5=1970", 2)
$d.Content.Find.Execute("983×2=1966", $true, $false, $false, $false, $false, $true, 1, $false, "867×6=5202", 2)
